$wb = $excel.ActiveWorkbook

# ---- Sheet: 大智投资组合收益 (rows 94-104) ----
$ws = $wb.Worksheets.Item(1)

$data = @(
    @('大智 (稳健智远)', '000333', '美的集团', 3.03, 42.89719527444501, 71.33, 3059.856938926162, 100887.2158440456, '202506191600'),
    @('大智 (稳健智远)', '510050', '上证50ETF', 5.01, 1852.638869852698, 2.73, 5057.704114697865, 100887.2158440456, '202506191600'),
    @('大智 (稳健智远)', '510300', '沪深300ETF', 4.9, 1276.881426590205, 3.87, 4941.531120904093, 100887.2158440456, '202506191600'),
    @('大智 (稳健智远)', '518880', '黄金ETF', 4.98, 673.9096418114974, 7.46, 5027.365927913771, 100887.2158440456, '202506191600'),
    @('大智 (稳健智远)', '600085', '同仁堂', 1.99, 56.06334956913254, 35.85, 2009.871082053402, 100887.2158440456, '202506191600'),
    @('大智 (稳健智远)', '600900', '长江电力', 20.37, 669.0422707938175, 30.72, 20552.97855878607, 100887.2158440456, '202506191600'),
    @('大智 (稳健智远)', '600989', '宝丰能源', 5.11, 308.5861230826723, 16.7, 5153.388255480628, 100887.2158440456, '202506191600'),
    @('大智 (稳健智远)', '601899', 'XD紫金矿', 9.859999999999999, 541.4194359293219, 18.37, 9945.875038021644, 100887.2158440456, '202506191600'),
    @('大智 (稳健智远)', 'HK02899', '紫金矿业', 9.82, 521.4694874201556, 19, 9907.920260982955, 100887.2158440456, '202506191600'),
    @('大智 (稳健智远)', 'HK06881', '中国银河', 4.62, 600.7968033130801, 7.76, 4662.183193709501, 100887.2158440456, '202506191600'),
    @('大智 (稳健智远)', '100000', '现金', 30.3, 30568.54135256952, 1, 30568.54135256952, 100887.2158440456, '202506191600')
)

$startRow = 94
$rowCount = $data.Count
$colCount = 9

# Force columns B and I to be stored as text (stock codes / timestamps keep leading zeros / aren't numbers)
$rngB = $ws.Range("B94:B104")
$rngI = $ws.Range("I94:I104")
$rngB.NumberFormat = "@"
$rngI.NumberFormat = "@"

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 0; $c -lt $colCount; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

# Reset style to Normal (default) so only the number format we need is applied, matching a plain unstyled cell
$rngB.Style = "Normal"
$rngI.Style = "Normal"

# ---- Sheet: 大成投资组合收益 (rows 65-71) ----
$ws = $wb.Worksheets.Item(2)

$data = @(
    @('大成 (锐进先锋)', '000725', '京东方A', 5.06, 1264.047242468624, 3.93, 4967.705662901692, 98175.3787773389, '202506191600'),
    @('大成 (锐进先锋)', '159781', '科创创业ETF', 5.01, 9277.629760760277, 0.53, 4917.143773202947, 98175.3787773389, '202506191600'),
    @('大成 (锐进先锋)', '513100', '纳指ETF', 4.98, 3131.938709046463, 1.56, 4885.824386112482, 98175.3787773389, '202506191600'),
    @('大成 (锐进先锋)', '513290', '纳指生物科技ETF', 0.98, 870.2909333102562, 1.1, 957.3200266412819, 98175.3787773389, '202506191600'),
    @('大成 (锐进先锋)', '603119', '浙江荣泰', 46.19, 1069.978093782073, 42.38, 45345.67161448426, 98175.3787773389, '202506191600'),
    @('大成 (锐进先锋)', '688290', '景业智能', 7.74, 147.121987895035, 51.65, 7598.850674778558, 98175.3787773389, '202506191600'),
    @('大成 (锐进先锋)', '100000', '现金', 30.05, 29502.86263921768, 1, 29502.86263921768, 98175.3787773389, '202506191600')
)

$startRow = 65
$rowCount = $data.Count
$colCount = 9

# Force columns B and I to be stored as text (stock codes / timestamps keep leading zeros / aren't numbers)
$rngB = $ws.Range("B65:B71")
$rngI = $ws.Range("I65:I71")
$rngB.NumberFormat = "@"
$rngI.NumberFormat = "@"

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 0; $c -lt $colCount; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

# Reset style to Normal (default) so only the number format we need is applied, matching a plain unstyled cell
$rngB.Style = "Normal"
$rngI.Style = "Normal"

# ---- Sheet: 我的投资组合收益 (rows 136-150) ----
$ws = $wb.Worksheets.Item(3)

$data = @(
    @('范式进化投资组合', '000333', '美的集团', 1.01, 14.02515882310653, 71.33, 1000.414578852189, 99516.0052923568, '202506191600'),
    @('范式进化投资组合', '000725', '京东方A', 5.07, 1284.618017653691, 3.93, 5048.548809379005, 99516.0052923568, '202506191600'),
    @('范式进化投资组合', '159781', '科创创业ETF', 5.02, 9428.611488061997, 0.53, 4997.164088672858, 99516.0052923568, '202506191600'),
    @('范式进化投资组合', '510050', '上证50ETF', 4.98, 1817.150577699221, 2.73, 4960.821077118873, 99516.0052923568, '202506191600'),
    @('范式进化投资组合', '510300', '沪深300ETF', 4.87, 1252.422077361618, 3.87, 4846.873439389462, 99516.0052923568, '202506191600'),
    @('范式进化投资组合', '513100', '纳指ETF', 1, 636.5814125697908, 1.56, 993.0670036088737, 99516.0052923568, '202506191600'),
    @('范式进化投资组合', '513290', '纳指生物科技ETF', 0.98, 884.4538210040457, 1.1, 972.8992031044504, 99516.0052923568, '202506191600'),
    @('范式进化投资组合', '518880', '黄金ETF', 0.99, 132.2001081659486, 7.46, 986.2128069179765, 99516.0052923568, '202506191600'),
    @('范式进化投资组合', '600085', '同仁堂', 0.99, 27.49471300507762, 35.85, 985.6854612320327, 99516.0052923568, '202506191600'),
    @('范式进化投资组合', '600900', '长江电力', 1.01, 32.81132034584936, 30.72, 1007.963761024492, 99516.0052923568, '202506191600'),
    @('范式进化投资组合', '600989', '宝丰能源', 5.08, 302.6749902285196, 16.7, 5054.672336816277, 99516.0052923568, '202506191600'),
    @('范式进化投资组合', '601899', 'XD紫金矿', 9.800000000000001, 531.0482559694855, 18.37, 9755.356462159449, 99516.0052923568, '202506191600'),
    @('范式进化投资组合', '603119', '浙江荣泰', 1.03, 24.16423640557475, 42.38, 1024.080338868258, 99516.0052923568, '202506191600'),
    @('范式进化投资组合', 'HK06881', '中国银河', 0.92, 117.8576436007749, 7.76, 914.5753143420133, 99516.0052923568, '202506191600'),
    @('范式进化投资组合', '100000', '现金', 57.24, 56967.67061087058, 1, 56967.67061087058, 99516.0052923568, '202506191600')
)

$startRow = 136
$rowCount = $data.Count
$colCount = 9

# Force columns B and I to be stored as text (stock codes / timestamps keep leading zeros / aren't numbers)
$rngB = $ws.Range("B136:B150")
$rngI = $ws.Range("I136:I150")
$rngB.NumberFormat = "@"
$rngI.NumberFormat = "@"

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 0; $c -lt $colCount; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

# Reset style to Normal (default) so only the number format we need is applied, matching a plain unstyled cell
$rngB.Style = "Normal"
$rngI.Style = "Normal"

